{"js": "const newValues = [\n  \"91-56=\",\n  \"58+6=\",\n  \"35+6=\",\n  \"30-8=\",\n  \"58+33=\",\n  \"44-7=\",\n  \"15+77=\",\n  \"81-76=\",\n  \"5+89=\",\n  \"13+28=\",\n  \"63-36=\",\n  \"65-38=\",\n  \"23+69=\",\n  \"50-1=\",\n  \"79+5=\",\n  \"46-19=\",\n  \"73-27=\",\n  \"53-19=\",\n  \"66+6=\",\n  \"91-73=\",\n  \"84-48=\",\n  \"49+29=\",\n  \"44-15=\",\n  \"4+58=\",\n  \"94-66=\",\n  \"6+35=\",\n  \"52-25=\",\n  \"96-89=\",\n  \"85-49=\",\n  \"29+22=\",\n  \"63-15=\",\n  \"55+6=\",\n  \"50-31=\",\n  \"47-8=\",\n  \"25+18=\",\n  \"58+35=\",\n  \"80-76=\",\n  \"47+16=\",\n  \"2+39=\",\n  \"53-26=\",\n  \"54-16=\",\n  \"38+23=\",\n  \"9+54=\",\n  \"26+57=\",\n  \"83-38=\",\n  \"72-36=\",\n  \"22+69=\",\n  \"54+8=\",\n  \"18+43=\",\n  \"75+9=\",\n  \"59+7=\",\n  \"48+9=\",\n  \"58+13=\",\n  \"71-48=\",\n  \"92-58=\",\n  \"13-9=\",\n  \"59+35=\",\n  \"23+69=\",\n  \"93-29=\",\n  \"19+25=\",\n  \"29+2=\",\n  \"32-19=\",\n  \"42-29=\",\n  \"67+24=\",\n  \"13+49=\",\n  \"49+23=\",\n  \"76-48=\",\n  \"72-9=\",\n  \"76-49=\",\n  \"95-58=\",\n  \"30-25=\",\n  \"15+8=\",\n  \"52+39=\",\n  \"8+69=\",\n  \"69+7=\",\n  \"11-3=\",\n  \"67+27=\",\n  \"16+75=\",\n  \"7+44=\",\n  \"18+76=\",\n  \"77-39=\",\n  \"81-53=\",\n  \"61-57=\",\n  \"64-15=\",\n  \"72-54=\",\n  \"8+19=\",\n  \"86-59=\",\n  \"34+38=\",\n  \"29+32=\",\n  \"18+43=\",\n  \"88-9=\",\n  \"42-18=\",\n  \"71-6=\",\n  \"48+14=\",\n  \"45-9=\",\n  \"61-9=\",\n  \"3+78=\",\n  \"38+4=\",\n  \"92-19=\",\n  \"9+17=\",\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst cols = 5;\nlet idx = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < cols; c++) {\n    const cell = table.getCell(r, c);\n    cell.getRange().insertText(newValues[idx], \"Replace\");\n    idx++;\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$newValues = @(\n  '91-56=',\n  '58+6=',\n  '35+6=',\n  '30-8=',\n  '58+33=',\n  '44-7=',\n  '15+77=',\n  '81-76=',\n  '5+89=',\n  '13+28=',\n  '63-36=',\n  '65-38=',\n  '23+69=',\n  '50-1=',\n  '79+5=',\n  '46-19=',\n  '73-27=',\n  '53-19=',\n  '66+6=',\n  '91-73=',\n  '84-48=',\n  '49+29=',\n  '44-15=',\n  '4+58=',\n  '94-66=',\n  '6+35=',\n  '52-25=',\n  '96-89=',\n  '85-49=',\n  '29+22=',\n  '63-15=',\n  '55+6=',\n  '50-31=',\n  '47-8=',\n  '25+18=',\n  '58+35=',\n  '80-76=',\n  '47+16=',\n  '2+39=',\n  '53-26=',\n  '54-16=',\n  '38+23=',\n  '9+54=',\n  '26+57=',\n  '83-38=',\n  '72-36=',\n  '22+69=',\n  '54+8=',\n  '18+43=',\n  '75+9=',\n  '59+7=',\n  '48+9=',\n  '58+13=',\n  '71-48=',\n  '92-58=',\n  '13-9=',\n  '59+35=',\n  '23+69=',\n  '93-29=',\n  '19+25=',\n  '29+2=',\n  '32-19=',\n  '42-29=',\n  '67+24=',\n  '13+49=',\n  '49+23=',\n  '76-48=',\n  '72-9=',\n  '76-49=',\n  '95-58=',\n  '30-25=',\n  '15+8=',\n  '52+39=',\n  '8+69=',\n  '69+7=',\n  '11-3=',\n  '67+27=',\n  '16+75=',\n  '7+44=',\n  '18+76=',\n  '77-39=',\n  '81-53=',\n  '61-57=',\n  '64-15=',\n  '72-54=',\n  '8+19=',\n  '86-59=',\n  '34+38=',\n  '29+32=',\n  '18+43=',\n  '88-9=',\n  '42-18=',\n  '71-6=',\n  '48+14=',\n  '45-9=',\n  '61-9=',\n  '3+78=',\n  '38+4=',\n  '92-19=',\n  '9+17='\n)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$idx]\n    $idx++\n  }\n}"}
